$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Arkusz1")

# New order rows added under the "Kosztorys" table (pinout + mainboard invoice items)
$ws.Cells.Item(15, 2).Value = "Nakrętki do profili"
$ws.Cells.Item(15, 3).Value = "AV-Elektronika"
$ws.Cells.Item(15, 4).Value = 24.5

$ws.Cells.Item(16, 2).Value = "Części do mainboard'a"
$ws.Cells.Item(16, 3).Value = "TME"
$ws.Cells.Item(16, 4).Value = 277.74

$ws.Range("E15").Select()
